$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")
$scratch.Font.Name = "Times New Roman"
$scratch.Font.Size = 10
$scratch.WrapText = $true
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$ws.Range("F165").PasteSpecial(-4122)
Write-Output "done"
